$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")

# --- Rename existing quarter headers from "Qn '2x" style to "Qn FY2x" style ---
$ws.Range("B1").Value = "Q1 FY23"
$ws.Range("C1").Value = "Q2 FY23"
$ws.Range("D1").Value = "Q3 FY23"
$ws.Range("E1").Value = "Q4 FY23"
$ws.Range("F1").Value = "Q1 FY24"
$ws.Range("G1").Value = "Q2 FY24"
$ws.Range("H1").Value = "Q3 FY24"
$ws.Range("I1").Value = "Q4 FY24"

# --- Add the three new reporting quarters (FY25 Q1-Q3) ---
$ws.Range("J1").Value = "Q1 FY25"
$ws.Range("K1").Value = "Q2 FY25"
$ws.Range("L1").Value = "Q3 FY25"

# --- Report Date row: copy the date formatting from the last existing quarter (I2) ---
$ws.Range("I2").Copy() | Out-Null
$ws.Range("J2:L2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("J2").Value = 45690
$ws.Range("K2").Value = 45781
$ws.Range("L2").Value = 45872

# --- Semiconductor Solutions row: copy the integer formatting from the last existing quarter (I3) ---
$ws.Range("I3").Copy() | Out-Null
$ws.Range("J3:L3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("J3").Value = 8212
$ws.Range("K3").Value = 8408
$ws.Range("L3").Value = 9166

# --- Infrastructure Software row: copy the integer formatting from the last existing quarter (I4) ---
$ws.Range("I4").Copy() | Out-Null
$ws.Range("J4:L4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("J4").Value = 6704
$ws.Range("K4").Value = 6596
$ws.Range("L4").Value = 6786

# --- Leave the selection/active cell on L4, matching the last edited cell ---
$ws.Range("L4").Select() | Out-Null
